# Updates Price (D) and Volume(1h) (E) columns for the cryptos list,
# reflecting refreshed market data from the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row|D-value|E-value  (D/E left blank means "no change for that column";
# a leading "'" forces text so Excel doesn't reinterpret numeric-looking
# price strings like "233.41" as a Number, matching the source data which
# stores these as plain text cells).
$data = @"
2|37.900.03|  +1.45%  
3|2.093.33|  +1.05%  
4||  +0.02%  
5|'233.41|  -0.17%  
6|'0.626|
7||  -0.03%  
8|'57.62|  +1.00%  
9||  +1.89%  
10|'0.0783|  +2.54%  
11||  +2.73%  
12|2.390.81|  +0.75%  
13|'14.43|  -1.45%  
14|'21.20|  +2.22%  
15|'0.763|  -1.85%  
16|'5.25|  +2.24%  
17|2.093.22|  +1.21%  
18|37.840.57|  +1.51%  
19||  -2.99%  
20||  +1.97%  
21|0.0₃0823|  +1.31%  
22|'228.60|  +1.01%  
23||  +0.02%  
24||  -0.97%  
25|'2.39|  -0.32%  
26|'170.58|  +1.75%  
27|'0.140|  +10.87%  
28|'8.94|  +1.94%  
29||  +0.13%  
30||  +2.19%  
31||  +1.37%  
32|'4.64|  +4.10%  
33|'0.0629|  +2.04%  
34|'4.60|  +0.81%  
35|'2.53|  +1.34%  
36|'1.83|  +4.02%  
37||  +4.99%  
38||  -0.06%  
39|'5.44|  -3.87%  
40||  +7.16%  
41|'2.94|  -0.48%  
42|'97.43|  +1.19%  
43|'0.0214|  +1.02%  
44|1.453.91|  -1.27%  
45||  -0.28%  
46||  +3.60%  
47|'15.78|  +4.67%  
48|'4.04|  -8.37%  
49|'7.41|  +3.71%  
50|'3.02|  +2.05%  
51|2.285.38|  +1.11%  
"@

$lines = $data -split "`n"
foreach ($rawLine in $lines) {
    $line = $rawLine.TrimEnd("`r")
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|', 3
    $row = [int]$parts[0]
    $dval = $parts[1]
    $eval = $parts[2]
    if ($dval.Length -gt 0) {
        $ws.Cells.Item($row, 4).Value = $dval
    }
    if ($eval.Length -gt 0) {
        $ws.Cells.Item($row, 5).Value = $eval
    }
}

Write-Output "Applied price/volume updates to $($lines.Count) rows"
